$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet ("tasks") and name it "documents"
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$docs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$docs.Name = "documents"

# Header row (yellow fill, same style used by the other sheets' header rows)
$docs.Range("A1:C1").Interior.ColorIndex = 6
$docs.Range("A1").Value = "title"
$docs.Range("B1").Value = "description"
$docs.Range("C1").Value = "version"

# Data rows - fill column by column to mirror how the data was authored
$docs.Range("A2").Value = "Title - 1"
$docs.Range("A3").Value = "Title - 2"

$docs.Range("B2").Value = "Desc - 1"
$docs.Range("B3").Value = "Desc - 2"

$docs.Range("C2").Value = "Version -3"
$docs.Range("C3").Value = "Version -4"

# Size the columns to fit their contents
[void]$docs.Columns.Item(1).AutoFit()
[void]$docs.Columns.Item(2).AutoFit()
$docs.Columns.Item(3).ColumnWidth = 12

# Make the new sheet the active/selected tab, matching the authored selection
[void]$docs.Activate()
[void]$docs.Range("H5").Select()
